# Updates column C ("Fitness") values in Sheet1 to reflect the new run data.
# Each contiguous block of rows is set via a single Range.Value assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2:C8").Value = 8693
$ws.Range("C9:C19").Value = 8486
$ws.Range("C20:C26").Value = 8484
$ws.Range("C27:C30").Value = 7919
$ws.Range("C31:C51").Value = 7917
$ws.Range("C52:C52").Value = 7912
$ws.Range("C53:C56").Value = 7828
$ws.Range("C57:C103").Value = 7639
$ws.Range("C104:C116").Value = 7312
$ws.Range("C117:C252").Value = 7310
